$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Duel Decks Elves vs. Inventors (DDU)',
    'Albino Troll',
    'Bad Moon',
    'Basking Rootwalla',
    'Beast Attack',
    'Blastoderm',
    'Corrupt',
    'Deathgreeter',
    'Drudge Skeletons',
    'Elephant Guide',
    'Enslave',
    'Faerie Macabre',
    'Fleshbag Marauder',
    'Forest',
    'Forest',
    'Forest',
    'Forest',
    'Garruk Wildspeaker',
    'Genju of the Cedars',
    'Genju of the Fens',
    'Ghost-Lit Stalker',
    'Giant Growth',
    'Harmonize',
    'Hideous End',
    'Howling Banshee',
    'Ichor Slick',
    'Indrik Stomphowler',
    'Invigorate',
    'Keening Banshee',
    'Krosan Tusker',
    'Lignify',
    'Liliana Vess',
    'Mutilate',
    'Nature''s Lore',
    'Overrun',
    'Phyrexian Rager',
    'Plated Slagwurm',
    'Polluted Mire',
    'Rancor',
    'Ravenous Baloth',
    'Ravenous Rats',
    'Rise from the Grave',
    'Rude Awakening',
    'Serrated Arrows',
    'Sign in Blood',
    'Skeletal Vampire',
    'Slippery Karst',
    'Snuff Out',
    'Stampeding Wildebeests',
    'Swamp',
    'Swamp',
    'Swamp',
    'Swamp',
    'Tendrils of Corruption',
    'Treetop Village',
    'Twisted Abomination',
    'Urborg Syphon-Mage',
    'Vampire Bats',
    'Vicious Hunger',
    'Vine Trellis',
    'Wall of Bone',
    'Wild Mongrel',
    'Windstorm',
    'Wirewood Savage'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
